# "last changes to v1.8.2"
#
# Bumps the StructureDefinition metadata to version 1.8.2 / new publish
# date, and fills in the previously-missing ele-1/ext-1 invariant text
# on the root Extension row of the Elements sheet.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.8.2"
$meta.Range("B8").Value = "2023-09-01T14:45:29-04:00"

# --- Elements sheet ---------------------------------------------------
# Row 1 is the root "Extension" element; its Invariants column (AJ) was
# blank and now carries the same ele-1 / ext-1 constraint text already
# present on the Extension.extension row.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AJ1").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
